$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.687.43"
$ws.Range("E2").Value = "  +4.47%  "

$ws.Range("D3").Value = "2.433.82"
$ws.Range("E3").Value = "  +6.27%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "481.06"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +8.68%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.93"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +15.02%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.23%  "

$ws.Range("E8").Value = "  +9.13%  "

$ws.Range("D9").Value = "2.457.36"
$ws.Range("E9").Value = "  +7.53%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0967"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +11.54%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.46"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +3.69%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.325"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +9.09%  "

$ws.Range("E13").Value = "  +2.02%  "

$ws.Range("D14").Value = "2.867.24"
$ws.Range("E14").Value = "  +6.35%  "

$ws.Range("D15").Value = "54.881.29"
$ws.Range("E15").Value = "  +4.77%  "

$ws.Range("E16").Value = "  +9.98%  "

$ws.Range("E17").Value = "  +14.82%  "

$ws.Range("D18").Value = "2.463.01"
$ws.Range("E18").Value = "  +5.85%  "

$ws.Range("E19").Value = "  +11.43%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "314.68"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +6.84%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.82"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +11.40%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.995"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.35%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.67"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +11.51%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "57.26"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +7.56%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.166"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +13.33%  "

$ws.Range("E26").Value = "  +12.70%  "

$ws.Range("E27").Value = "  +0.56%  "

$ws.Range("D28").Value = "2.596.77"
$ws.Range("E28").Value = "  +7.71%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.30"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +7.69%  "

$ws.Range("D30").Value = "0.0₃0777"
$ws.Range("E30").Value = "  +17.76%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.998"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.09%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "148.12"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +5.13%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "17.95"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +7.45%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.47"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +10.88%  "

$ws.Range("E35").Value = "  +9.36%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.10"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +12.55%  "

$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.63"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +7.53%  "

$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.854"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +4.44%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "33.80"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +7.24%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.993"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.13%  "

$ws.Range("E41").Value = "  +6.88%  "

$ws.Range("E42").Value = "  +9.75%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.39"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +9.74%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.29"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +10.08%  "

$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.67"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +12.44%  "

$ws.Range("B46").Value = "WhiteBITCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.09"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.52%  "

$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0903"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +10.73%  "

$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "1.952.99"
$ws.Range("E48").Value = "  +3.36%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "252.92"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +30.70%  "

$ws.Range("E50").Value = "  +8.58%  "

$ws.Range("E51").Value = "  +10.94%  "
